$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 261386.58
$ws.Range("I98").Value = 339833.53
$ws.Range("J98").Value = 2511.6
$ws.Range("K98").Value = 339833.53
$ws.Range("L98").Value = 2511.6
$ws.Range("M98").Value = -338335.53
$ws.Range("N98").Value = -5507.6
$ws.Range("H112").Value = 6945425.5
$ws.Range("I112").Value = 866.6667
$ws.Range("J112").Value = 7576749
$ws.Range("K112").Value = 2600.0001
$ws.Range("L112").Value = 22730247
$ws.Range("M112").Value = -1492.0001
$ws.Range("N112").Value = -22732463
$ws.Range("H122").Value = 261386.58
$ws.Range("I122").Value = 339833.53
$ws.Range("J122").Value = 2511.6
$ws.Range("K122").Value = 1019500.59
$ws.Range("L122").Value = 7534.799999999999
$ws.Range("M122").Value = -1017050.59
$ws.Range("N122").Value = -12434.8
$ws.Range("H128").Value = 79660
$ws.Range("J128").Value = 79660
$ws.Range("L128").Value = 79660
$ws.Range("N128").Value = -89620
$ws.Range("H129").Value = 1157.9565
$ws.Range("I129").Value = 536.8
$ws.Range("J129").Value = 1330.5
$ws.Range("K129").Value = 1610.4
$ws.Range("L129").Value = 3991.5
$ws.Range("M129").Value = 3389.6
$ws.Range("N129").Value = -13991.5
$ws.Range("I137").Value = 76924320
$ws.Range("J137").Value = 1833.3334
$ws.Range("K137").Value = 230772960
$ws.Range("L137").Value = 5500.0002
$ws.Range("M137").Value = -230770410
$ws.Range("N137").Value = -10600.0002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46817.773
$ws.Range("I2").Value = 53852.156
$ws.Range("J2").Value = 2266.6667
$ws.Range("K2").Value = 53852.156
$ws.Range("L2").Value = 2266.6667
$ws.Range("M2").Value = -53739.156
$ws.Range("N2").Value = -2492.6667
$ws.Range("H32").Value = 13599.738
$ws.Range("I32").Value = 2292.9614
$ws.Range("K32").Value = 2292.9614
$ws.Range("M32").Value = -2005.9614
$ws.Range("H45").Value = 934
$ws.Range("I45").Value = 827.4
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 827.4
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -450.4
$ws.Range("N45").Value = -2754
$ws.Range("H74").Value = 13868.667
$ws.Range("I74").Value = 1852.25
$ws.Range("J74").Value = 110000
$ws.Range("K74").Value = 1852.25
$ws.Range("L74").Value = 110000
$ws.Range("M74").Value = -978.25
$ws.Range("N74").Value = -111748
$ws.Range("H77").Value = 13868.667
$ws.Range("I77").Value = 1852.25
$ws.Range("J77").Value = 110000
$ws.Range("K77").Value = 9261.25
$ws.Range("L77").Value = 550000
$ws.Range("M77").Value = -4893.25
$ws.Range("N77").Value = -558736
$ws.Range("H116").Value = 46817.773
$ws.Range("I116").Value = 53852.156
$ws.Range("J116").Value = 2266.6667
$ws.Range("K116").Value = 53852.156
$ws.Range("L116").Value = 2266.6667
$ws.Range("M116").Value = -51558.156
$ws.Range("N116").Value = -6854.6667
$ws.Range("H122").Value = 2309.524
$ws.Range("I122").Value = 2159
$ws.Range("J122").Value = 2402.1538
$ws.Range("K122").Value = 6477
$ws.Range("L122").Value = 7206.4614
$ws.Range("M122").Value = -4027
$ws.Range("N122").Value = -12106.4614

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 46817.773
$ws.Range("I3").Value = 53852.156
$ws.Range("J3").Value = 2266.6667
$ws.Range("K3").Value = 53852.156
$ws.Range("L3").Value = 2266.6667
$ws.Range("M3").Value = -53738.156
$ws.Range("N3").Value = -2494.6667
$ws.Range("H94").Value = 1587.2858
$ws.Range("I94").Value = 1472.2727
$ws.Range("J94").Value = 2009
$ws.Range("K94").Value = 1472.2727
$ws.Range("L94").Value = 2009
$ws.Range("M94").Value = -1021.2727
$ws.Range("N94").Value = -2911

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 49500
$ws.Range("J14").Value = 49500
$ws.Range("L14").Value = 49500
$ws.Range("N14").Value = -49840
$ws.Range("H58").Value = 982.78845
$ws.Range("I58").Value = 601.64514
$ws.Range("J58").Value = 1545.4286
$ws.Range("K58").Value = 601.64514
$ws.Range("L58").Value = 1545.4286
$ws.Range("M58").Value = -398.64514
$ws.Range("N58").Value = -1951.4286
$ws.Range("H107").Value = 331.75
$ws.Range("I107").Value = 198
$ws.Range("J107").Value = 376.33334
$ws.Range("K107").Value = 198
$ws.Range("L107").Value = 376.33334
$ws.Range("M107").Value = 1722
$ws.Range("N107").Value = -4216.33334
$ws.Range("H136").Value = 982.78845
$ws.Range("I136").Value = 601.64514
$ws.Range("J136").Value = 1545.4286
$ws.Range("K136").Value = 1804.93542
$ws.Range("L136").Value = 4636.2858
$ws.Range("M136").Value = 745.0645800000002
$ws.Range("N136").Value = -9736.2858
$ws.Range("H137").Value = 49945
$ws.Range("J137").Value = 49945
$ws.Range("L137").Value = 49945
$ws.Range("N137").Value = -60145
$ws.Range("H140").Value = 55981.668
$ws.Range("J140").Value = 55981.668
$ws.Range("L140").Value = 55981.668
$ws.Range("N140").Value = -66341.66800000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7038.6665
$ws.Range("I70").Value = 8373.75
$ws.Range("K70").Value = 8373.75
$ws.Range("M70").Value = -8103.75
$ws.Range("H73").Value = 7038.6665
$ws.Range("I73").Value = 8373.75
$ws.Range("K73").Value = 8373.75
$ws.Range("M73").Value = -7437.75
$ws.Range("H80").Value = 2411.56
$ws.Range("I80").Value = 2414.5715
$ws.Range("J80").Value = 2407.7273
$ws.Range("K80").Value = 2414.5715
$ws.Range("L80").Value = 2407.7273
$ws.Range("M80").Value = -1416.5715
$ws.Range("N80").Value = -4403.7273
$ws.Range("H83").Value = 2411.56
$ws.Range("I83").Value = 2414.5715
$ws.Range("J83").Value = 2407.7273
$ws.Range("K83").Value = 12072.8575
$ws.Range("L83").Value = 12038.6365
$ws.Range("M83").Value = -7080.8575
$ws.Range("N83").Value = -22022.6365
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H113").Value = 1786.3684
$ws.Range("I113").Value = 1835.2727
$ws.Range("J113").Value = 1719.125
$ws.Range("K113").Value = 1835.2727
$ws.Range("L113").Value = 1719.125
$ws.Range("M113").Value = 334.7273
$ws.Range("N113").Value = -6059.125
$ws.Range("H122").Value = 1113090.9
$ws.Range("I122").Value = 3704703.8
$ws.Range("J122").Value = 2399.7144
$ws.Range("K122").Value = 11114111.4
$ws.Range("L122").Value = 7199.1432
$ws.Range("M122").Value = -11111661.4
$ws.Range("N122").Value = -12099.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5070.3213
$ws.Range("I61").Value = 5742.95
$ws.Range("J61").Value = 3388.75
$ws.Range("K61").Value = 5742.95
$ws.Range("L61").Value = 3388.75
$ws.Range("M61").Value = -5540.95
$ws.Range("N61").Value = -3792.75
$ws.Range("H101").Value = 17633.5
$ws.Range("J101").Value = 17633.5
$ws.Range("L101").Value = 17633.5
$ws.Range("N101").Value = -24123.5
$ws.Range("H113").Value = 5070.3213
$ws.Range("I113").Value = 5742.95
$ws.Range("J113").Value = 3388.75
$ws.Range("K113").Value = 5742.95
$ws.Range("L113").Value = 3388.75
$ws.Range("M113").Value = -3572.95
$ws.Range("N113").Value = -7728.75
$ws.Range("H122").Value = 3357.1155
$ws.Range("I122").Value = 2342.8572
$ws.Range("J122").Value = 3730.7896
$ws.Range("K122").Value = 7028.571599999999
$ws.Range("L122").Value = 11192.3688
$ws.Range("M122").Value = -4578.571599999999
$ws.Range("N122").Value = -16092.3688

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 79923.30499999999
$ws.Range("I122").Value = 93091.17999999999
$ws.Range("K122").Value = 279273.54
$ws.Range("M122").Value = -276823.54
$ws.Range("H126").Value = 86508.414
$ws.Range("I126").Value = 103080.1
$ws.Range("J126").Value = 3650
$ws.Range("K126").Value = 309240.3
$ws.Range("L126").Value = 10950
$ws.Range("M126").Value = -306770.3
$ws.Range("N126").Value = -15890
